# Update "想去人数" (number of people interested) counts on two sheets.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (first sheet): rows 3-5, column F
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 177
$wsExhibit.Range("F4").Value = 769
$wsExhibit.Range("F5").Value = 66

# Sheet "全部类型" (fourth sheet): rows 4-6, column F
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 177
$wsAll.Range("F5").Value = 769
$wsAll.Range("F6").Value = 66
